$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy header style (bold, style index 1) from existing header cells to the new section header rows
$ws.Range("A125").Copy() | Out-Null
$ws.Range("A151").PasteSpecial(-4122) | Out-Null
$ws.Range("B125").Copy() | Out-Null
$ws.Range("B151").PasteSpecial(-4122) | Out-Null
$ws.Range("B126").Copy() | Out-Null
$ws.Range("B152").PasteSpecial(-4122) | Out-Null
$ws.Range("C126").Copy() | Out-Null
$ws.Range("C152").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the values for the new "Roles y Permisos" test-plan section
$ws.Range("A151").Value = 'Modulo'
$ws.Range("B151").Value = 'Roles y Permisos'
$ws.Range("B152").Value = 'Tests'
$ws.Range("C152").Value = 'Resultado Esperado'
$ws.Range("D152").Value = 'Resultado Otenido'
$ws.Range("A153").Value = 1
$ws.Range("B153").Value = 'Cada rol accede a modulos especificos a los que tiene permiso'
$ws.Range("C153").Value = 'un usuario, de ninguna forma debe obtener acceso un modulo'
$ws.Range("C154").Value = ' que no esta asignado a su rol o a sus roles.'
$ws.Range("A156").Value = 2
$ws.Range("B156").Value = 'Permisos de Administrador: '
$ws.Range("C156").Value = 'El administrador puede acceder a todos los modulos del sistema'
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = 'Permisos de Secreataria: '
$ws.Range("C158").Value = 'Que pueda acceder a cada uno de estos modulos y ningun otro.'
$ws.Range("B159").Value = '           can :manage, Area'
$ws.Range("B160").Value = '           can :manage, CargoManifest'
$ws.Range("B161").Value = '           can :anage, CargoManifestDetail '
$ws.Range("B162").Value = '           can :manage, City '
$ws.Range("B163").Value = '           can :manage, Country '
$ws.Range("B164").Value = '           can :manage, Customer '
$ws.Range("B165").Value = '           can :manage, CustomerType '
$ws.Range("B166").Value = '           can :manage, PaymentMethod'
$ws.Range("B167").Value = '           can :manage, Product'
$ws.Range("B168").Value = '           can :manage, ProductState'
$ws.Range("B169").Value = '           can :manage, ProductType'
$ws.Range("B170").Value = '           can :manage, Province'
$ws.Range("B171").Value = '           can :manage, Reason'
$ws.Range("B172").Value = '           can :manage, Receiver'
$ws.Range("B173").Value = '           can :manage, ReceiverAddress'
$ws.Range("B174").Value = '           can :manage, RetireNote '
$ws.Range("B175").Value = '           can :manage, RetireNoteState '
$ws.Range("B176").Value = '           can :manage, RoutingSheet '
$ws.Range("B177").Value = '           can :manage, RoutingSheetDetail '
$ws.Range("B178").Value = '           can :manage, RoutingSheetState '
$ws.Range("B179").Value = '           can :manage, ServiceType '
$ws.Range("B180").Value = '           can :manage, TransportGuide '
$ws.Range("B181").Value = '           can :manage, TransportGuideDetail '
$ws.Range("B182").Value = '           can :manage, TransportGuideState'
$ws.Range("A184").Value = 4
$ws.Range("B184").Value = 'Permisos de Entregador:'
$ws.Range("C184").Value = 'Que pueda acceder a cada uno de estos modulos y ningun otro.'
$ws.Range("B185").Value = '            can :manage, RoutingSheet'
$ws.Range("B186").Value = '            can :manage, RoutingSheetDetail'
$ws.Range("B187").Value = '            can :manage, RoutingSheetState'
$ws.Range("A189").Value = 5
$ws.Range("B189").Value = 'Timeout'
$ws.Range("C189").Value = 'Cada 10 minutos debe hacer un timeout, debe pedir al usuario que se loguee'
$ws.Range("C190").Value = 'luego deve volver a la misma pagina donde se encontraba cuando'
$ws.Range("C191").Value = 'se expiro la session.'

# Restore selection to match the final edited cell
$ws.Range("B151").Select() | Out-Null
